# This script re-shuffles the per-observation data (species id/name/author,
# red-list status, taxon id, count/unit/age-stage fields, locality name,
# coordinates and the public-comment field) among rows 3-16 of the active
# sheet. Each target row receives the full set of these fields from another
# (source) row, per the mapping below - derived by diffing the workbook's
# "before" and "after" OOXML. Columns not listed (C, M, N, O, S..AB, AD..AY)
# are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (the target row ends up holding what used to be
# the source row's data for the columns below)
$mapping = @{
    3  = 5
    4  = 10
    5  = 12
    6  = 14
    7  = 16
    8  = 3
    9  = 4
    10 = 6
    11 = 7
    12 = 8
    13 = 9
    14 = 11
    15 = 13
    16 = 15
}

$cols = @("A","B","D","E","F","G","H","I","J","K","L","P","Q","R","AC")

# Snapshot the current ("before") value of every relevant cell first, since
# the mapping below permutes rows in place and several target/source rows
# overlap (e.g. row 3 feeds row 8, row 5 feeds row 3, etc.) - writing as we
# go would clobber a value before it has been read for a later row.
$snapshot = @{}
foreach ($r in $mapping.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $srcVals[$c]
    }
}

Write-Output "done"
